# Updated symbol list on Sun Jan 29 03:51:54 UTC 2023 with GitHub Actions
#
# Refreshes the Price (column D) and Volume(1h) (column E) quotes for the
# crypto rows on Sheet1 to the latest scraped values. Values are written with a
# leading apostrophe so Excel stores them as literal text (matching the sheet's
# existing text-formatted numbers/percentages) instead of auto-converting them
# to numeric/percentage cell values, which would silently drop things like
# trailing zeros (e.g. "0.08100") or change the "--"/"--%" placeholder rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'307.95"
$ws.Range("E2").Value = "'-0.07%"

# Row 3
$ws.Range("D3").Value = "'39.69"
$ws.Range("E3").Value = "'1.39%"

# Row 4
$ws.Range("D4").Value = "'5.148"
$ws.Range("E4").Value = "'1.05%"

# Row 5
$ws.Range("D5").Value = "'0.08100"

# Row 6
$ws.Range("D6").Value = "'1.942"
$ws.Range("E6").Value = "'-2.79%"

# Row 7
$ws.Range("D7").Value = "'8.156"
$ws.Range("E7").Value = "'3.25%"

# Row 8
$ws.Range("D8").Value = "'4.227"
$ws.Range("E8").Value = "'1.25%"

# Row 9
$ws.Range("D9").Value = "'0.9283"
$ws.Range("E9").Value = "'-0.13%"

# Row 10
$ws.Range("D10").Value = "'0.1402"
$ws.Range("E10").Value = "'-0.90%"

# Row 11
$ws.Range("D11").Value = "'0.1918"
$ws.Range("E11").Value = "'-1.16%"

# Row 12
$ws.Range("D12").Value = "'0.09099"
$ws.Range("E12").Value = "'-1.54%"

# Row 13
$ws.Range("D13").Value = "'0.03493"
$ws.Range("E13").Value = "'0.34%"

# Row 14
$ws.Range("D14").Value = "'0.09809"
$ws.Range("E14").Value = "'-0.23%"

# Row 15
$ws.Range("D15").Value = "'0.001390"
$ws.Range("E15").Value = "'-1.10%"

# Row 16
$ws.Range("D16").Value = "'0.005958"
$ws.Range("E16").Value = "'0.14%"

# Row 17
$ws.Range("D17").Value = "'3.943"
$ws.Range("E17").Value = "'-0.12%"

# Row 19
$ws.Range("D19").Value = "'0.3428"
$ws.Range("E19").Value = "'-0.73%"

# Row 20
$ws.Range("D20").Value = "'0.1321"
$ws.Range("E20").Value = "'1.47%"

# Row 21
$ws.Range("D21").Value = "'4.666"
$ws.Range("E21").Value = "'-2.79%"

# Row 23
$ws.Range("D23").Value = "'0.04389"
$ws.Range("E23").Value = "'-1.80%"

# Row 24
$ws.Range("D24").Value = "'0.001232"
$ws.Range("E24").Value = "'-0.65%"

# Row 25
$ws.Range("D25").Value = "'0.004349"
$ws.Range("E25").Value = "'4.21%"

# Row 26
$ws.Range("D26").Value = "'0.0001301"
$ws.Range("E26").Value = "'0.03%"

# Row 27
$ws.Range("D27").Value = "'0.0004002"
$ws.Range("E27").Value = "'-10.02%"

# Row 39
$ws.Range("D39").Value = "'0.02030"
$ws.Range("E39").Value = "'-3.96%"

# Row 40
$ws.Range("D40").Value = "'0.05053"
$ws.Range("E40").Value = "'-2.15%"

# Row 41
$ws.Range("D41").Value = "'0.007357"
$ws.Range("E41").Value = "'-1.48%"

# Row 42
$ws.Range("D42").Value = "'0.009770"
$ws.Range("E42").Value = "'-3.69%"

# Row 43
$ws.Range("E43").Value = "'-0.32%"

# Row 44
$ws.Range("D44").Value = "'0.002132"
$ws.Range("E44").Value = "'0.03%"

# Row 45
$ws.Range("D45").Value = "'0.008686"
$ws.Range("E45").Value = "'-10.17%"

# Row 46
$ws.Range("D46").Value = "'0.00006350"
$ws.Range("E46").Value = "'0.77%"

# Row 47
$ws.Range("D47").Value = "'0.00000000750"
$ws.Range("E47").Value = "'0.06%"

# Row 48
$ws.Range("D48").Value = "'0.002864"

# Row 49
$ws.Range("E49").Value = "'-18.72%"

# Row 50
$ws.Range("D50").Value = "'0.00002101"
$ws.Range("E50").Value = "'0.06%"

# Row 51
$ws.Range("D51").Value = "'0.0002001"
$ws.Range("E51").Value = "'0.06%"
